# Added validation set method for lda / qda
# Updates the "Test errors" summary row and the LDA/QDA confusion matrices
# to reflect the new validation-set-based error estimates.

$wb = $excel.ActiveWorkbook

# --- Sheet: "Test errors" ---
$wsErr = $wb.Worksheets.Item("Test errors")
$wsErr.Range("C2").Value = 0.06190560981646955
$wsErr.Range("D2").Value = 0.5416666666666666
$wsErr.Range("E2").Value = 0.13350683156331017

# --- Sheet: "LDA Conf. Mat." ---
$wsLda = $wb.Worksheets.Item("LDA Conf. Mat.")
$wsLda.Range("B2").Value = 30.0
$wsLda.Range("C2").Value = 0.0
$wsLda.Range("D2").Value = 0.0
$wsLda.Range("E2").Value = 0.0
$wsLda.Range("F2").Value = 0.0
$wsLda.Range("G2").Value = 6.0

$wsLda.Range("B3").Value = 0.0
$wsLda.Range("C3").Value = 35.0
$wsLda.Range("D3").Value = 0.0
$wsLda.Range("E3").Value = 0.0
$wsLda.Range("F3").Value = 0.0
$wsLda.Range("G3").Value = 1.0

$wsLda.Range("B4").Value = 0.0
$wsLda.Range("C4").Value = 0.0
$wsLda.Range("D4").Value = 32.0
$wsLda.Range("E4").Value = 1.0
$wsLda.Range("F4").Value = 3.0
$wsLda.Range("G4").Value = 0.0

$wsLda.Range("B5").Value = 0.0
$wsLda.Range("C5").Value = 0.0
$wsLda.Range("D5").Value = 0.0
$wsLda.Range("E5").Value = 34.0
$wsLda.Range("F5").Value = 2.0
$wsLda.Range("G5").Value = 0.0

$wsLda.Range("B6").Value = 0.0
$wsLda.Range("C6").Value = 0.0
$wsLda.Range("D6").Value = 3.0
$wsLda.Range("E6").Value = 4.0
$wsLda.Range("F6").Value = 29.0
$wsLda.Range("G6").Value = 0.0

$wsLda.Range("B7").Value = 6.0
$wsLda.Range("C7").Value = 0.0
$wsLda.Range("D7").Value = 1.0
$wsLda.Range("E7").Value = 0.0
$wsLda.Range("F7").Value = 2.0
$wsLda.Range("G7").Value = 27.0

# --- Sheet: "QDA Conf. Mat." ---
$wsQda = $wb.Worksheets.Item("QDA Conf. Mat.")
$wsQda.Range("B2").Value = 12.0
$wsQda.Range("C2").Value = 0.0
$wsQda.Range("D2").Value = 0.0
$wsQda.Range("E2").Value = 0.0
$wsQda.Range("F2").Value = 0.0
$wsQda.Range("G2").Value = 24.0

$wsQda.Range("B3").Value = 0.0
$wsQda.Range("C3").Value = 20.0
$wsQda.Range("D3").Value = 5.0
$wsQda.Range("E3").Value = 4.0
$wsQda.Range("F3").Value = 0.0
$wsQda.Range("G3").Value = 7.0

$wsQda.Range("B4").Value = 0.0
$wsQda.Range("C4").Value = 2.0
$wsQda.Range("D4").Value = 18.0
$wsQda.Range("E4").Value = 3.0
$wsQda.Range("F4").Value = 7.0
$wsQda.Range("G4").Value = 6.0

$wsQda.Range("B5").Value = 0.0
$wsQda.Range("C5").Value = 1.0
$wsQda.Range("D5").Value = 5.0
$wsQda.Range("E5").Value = 15.0
$wsQda.Range("F5").Value = 13.0
$wsQda.Range("G5").Value = 2.0

$wsQda.Range("B6").Value = 0.0
$wsQda.Range("C6").Value = 1.0
$wsQda.Range("D6").Value = 5.0
$wsQda.Range("E6").Value = 11.0
$wsQda.Range("F6").Value = 13.0
$wsQda.Range("G6").Value = 6.0

$wsQda.Range("B7").Value = 6.0
$wsQda.Range("C7").Value = 3.0
$wsQda.Range("D7").Value = 2.0
$wsQda.Range("E7").Value = 1.0
$wsQda.Range("F7").Value = 3.0
$wsQda.Range("G7").Value = 21.0
